# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "71.148.09"
$ws.Range("E2").Value = "  +4.06%  "

# Row 3
$ws.Range("D3").Value = "2.625.74"
$ws.Range("E3").Value = "  +4.48%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.94"
$ws.Range("E5").Value = "  +2.43%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.17"
$ws.Range("E6").Value = "  +2.62%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("E8").Value = "  +1.43%  "

# Row 9
$ws.Range("D9").Value = "2.625.56"
$ws.Range("E9").Value = "  +4.49%  "

# Row 10
$ws.Range("E10").Value = "  +15.14%  "

# Row 11
$ws.Range("E11").Value = "  +0.07%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.347"
$ws.Range("E12").Value = "  +2.66%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.02"
$ws.Range("E13").Value = "  +0.63%  "

# Row 14
$ws.Range("D14").Value = "3.081.86"
$ws.Range("E14").Value = "  +3.78%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.63"
$ws.Range("E15").Value = "  +2.94%  "

# Row 16
$ws.Range("E16").Value = "  +7.41%  "

# Row 17
$ws.Range("D17").Value = "71.183.55"
$ws.Range("E17").Value = "  +4.51%  "

# Row 18
$ws.Range("D18").Value = "2.621.60"
$ws.Range("E18").Value = "  +4.37%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "382.01"
$ws.Range("E19").Value = "  +8.79%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.92"
$ws.Range("E20").Value = "  +5.53%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.54"
$ws.Range("E21").Value = "  +4.69%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.14"
$ws.Range("E22").Value = "  -1.05%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.45"
$ws.Range("E23").Value = "  +1.53%  "

# Row 24
$ws.Range("E24").Value = "  +5.55%  "

# Row 25
$ws.Range("E25").Value = "  +0.10%  "

# Row 26
$ws.Range("E26").Value = "  +8.01%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.61"
$ws.Range("E27").Value = "  +3.82%  "

# Row 28
$ws.Range("D28").Value = "2.760.68"
$ws.Range("E28").Value = "  +4.60%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.03%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0957"
$ws.Range("E30").Value = "  +6.35%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "541.44"
$ws.Range("E31").Value = "  +5.94%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.05"
$ws.Range("E32").Value = "  +2.78%  "

# Row 33
$ws.Range("E33").Value = "  +4.31%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.83"
$ws.Range("E34").Value = "  +2.86%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.01%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "165.49"
$ws.Range("E36").Value = "  +2.38%  "

# Row 37
$ws.Range("E37").Value = "  -1.74%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.22"
$ws.Range("E38").Value = "  +4.51%  "

# Row 39
$ws.Range("E39").Value = "  +7.22%  "

# Row 40
$ws.Range("E40").Value = "  +1.39%  "

# Row 41
$ws.Range("E41").Value = "  +3.91%  "

# Row 42
$ws.Range("E42").Value = "  +9.07%  "

# Row 43
$ws.Range("E43").Value = "  +0.11%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.05"
$ws.Range("E44").Value = "  +4.40%  "

# Row 45
$ws.Range("E45").Value = "  +0.76%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.09"
$ws.Range("E46").Value = "  +2.74%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "154.58"
$ws.Range("E47").Value = "  +2.39%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.64"
$ws.Range("E48").Value = "  +2.04%  "

# Row 49
$ws.Range("E49").Value = "  +1.91%  "

# Row 50
$ws.Range("E50").Value = "  +2.42%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.68"
$ws.Range("E51").Value = "  +5.20%  "
